$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: end year value was stored as a (wrong) negative number -2021; correct it to 2021
$ws.Range("D2").Value = 2021

# E2: end year for this grant is now a text note ",2022" instead of a plain number 2022
$ws.Range("E2").Value = ",2022"

# Update the active selection to E2 (matches the author's recorded cursor position)
$ws.Range("E2").Select()
